$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 19 ---
$ws.Range("M19").Value = 2.91371
$ws.Range("N19").Value = 1.7809200000000001
$ws.Range("U19").Value = 2.7817500000000002
$ws.Range("Y19").Value = 7.1057199999999998

# --- Row 20 ---
$ws.Range("M20").Value = 0.10781
$ws.Range("N20").Value = 0.14085
$ws.Range("U20").Value = 0.036
$ws.Range("Y20").Value = 0.501

# --- Row 29 ---
$ws.Range("N29").Value = 2.1948500000000002
$ws.Range("V29").Value = 2.2303099999999998

# --- Row 30 ---
$ws.Range("N30").Value = 0.046
$ws.Range("V30").Value = 0.15090000000000001

# --- Row 31 ---
$ws.Range("N31").Value = -0.13697000000000001
$ws.Range("Q31").Formula = "=N31*-1*1000"
$ws.Range("V31").Value = -0.15629999999999999
$ws.Range("Y31").Formula = "=V31*-1*1000"

# --- Row 32 ---
$ws.Range("N32").Value = 0.005
$ws.Range("Q32").NumberFormat = "0.00"
$ws.Range("Q32").Formula = "=N32*1000"
$ws.Range("V32").Value = 0.01841
$ws.Range("Y32").Formula = "=V32*1000"

# --- Selection / active cell ---
$ws.Range("S33").Select()
